$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    192.06396484375,
    190.4872436523438,
    185.0630950927734,
    186.1680908203125,
    183.5380554199219,
    177.8310546875,
    179.5807495117188,
    185.1297760009766,
    195.8791961669922,
    212.3919525146484,
    214.0627136230469,
    246.5023345947266,
    258.9142456054688,
    255.2961578369141,
    224.3261871337891,
    205.7127838134766,
    191.6334838867188,
    179.4331512451172,
    167.5816040039062,
    151.3641662597656,
    152.9140014648438,
    137.9481811523438,
    136.4194030761719,
    125.0732574462891
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
